# Regenerate save_data: recompute column G ("K") values for rows 2-35.
# (per commit: "regen save_data to use K instead of Strike#, regen std/mean,
#  calc and write s_vals")

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values for column G ("K") keyed by row number.
$kValues = @{
    2  = 1
    3  = 1
    4  = 2
    5  = 1
    6  = 0
    7  = 0
    8  = 0
    9  = 0
    10 = 0
    11 = 1
    12 = 0
    13 = 1
    14 = 2
    15 = 0
    16 = 1
    17 = 2
    18 = 2
    19 = 0
    20 = 1
    21 = 1
    22 = 1
    23 = 2
    24 = 1
    25 = 1
    26 = 2
    27 = 0
    28 = 2
    29 = 0
    30 = 1
    31 = 0
    32 = 1
    33 = 0
    34 = 0
    35 = 1
}

foreach ($row in $kValues.Keys) {
    $ws.Cells.Item($row, 7).Value = $kValues[$row]
}
